$wb = $excel.ActiveWorkbook

# --- Sheet "snapshot": update scraped_at (column K) timestamps for rows 2-34 ---
$snapshot = $wb.Worksheets.Item("snapshot")

$timestamps = @{
    2  = "2025-11-28T07:01:46.444715+00:00"
    3  = "2025-11-28T07:01:48.828176+00:00"
    4  = "2025-11-28T07:01:48.828209+00:00"
    5  = "2025-11-28T07:01:48.828227+00:00"
    6  = "2025-11-28T07:01:51.584186+00:00"
    7  = "2025-11-28T07:01:54.856430+00:00"
    8  = "2025-11-28T07:01:57.217986+00:00"
    9  = "2025-11-28T07:01:59.486195+00:00"
    10 = "2025-11-28T07:02:04.051106+00:00"
    11 = "2025-11-28T07:02:04.051148+00:00"
    12 = "2025-11-28T07:02:06.532836+00:00"
    13 = "2025-11-28T07:02:08.905927+00:00"
    14 = "2025-11-28T07:02:11.605414+00:00"
    15 = "2025-11-28T07:02:13.908911+00:00"
    16 = "2025-11-28T07:02:13.908939+00:00"
    17 = "2025-11-28T07:02:13.908958+00:00"
    18 = "2025-11-28T07:02:16.688791+00:00"
    19 = "2025-11-28T07:02:16.688822+00:00"
    20 = "2025-11-28T07:02:16.688840+00:00"
    21 = "2025-11-28T07:02:16.688856+00:00"
    22 = "2025-11-28T07:02:19.522043+00:00"
    23 = "2025-11-28T07:02:19.522072+00:00"
    24 = "2025-11-28T07:02:21.886861+00:00"
    25 = "2025-11-28T07:02:21.886892+00:00"
    26 = "2025-11-28T07:02:21.886909+00:00"
    27 = "2025-11-28T07:02:21.886928+00:00"
    28 = "2025-11-28T07:02:24.181498+00:00"
    29 = "2025-11-28T07:02:26.432659+00:00"
    30 = "2025-11-28T07:02:28.699798+00:00"
    31 = "2025-11-28T07:02:34.188496+00:00"
    32 = "2025-11-28T07:02:34.188527+00:00"
    33 = "2025-11-28T07:02:36.560639+00:00"
    34 = "2025-11-28T07:02:36.560674+00:00"
}

foreach ($row in $timestamps.Keys) {
    $snapshot.Cells.Item($row, 11).Value = $timestamps[$row]
}

# --- Sheet "new_injured": remove data rows 2-4, keep only header row ---
$newInjured = $wb.Worksheets.Item("new_injured")
$newInjured.Range("A2:G4").Delete()
